$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: headers - swap columns H and I (ExtractionType now before SamplePortion)
$ws.Cells.Item(1, 8).Value = "ExtractionType"
$ws.Cells.Item(1, 9).Value = "SamplePortion"

# Row 2: type hints - swap columns H and I, update float hint with unit
$ws.Cells.Item(2, 8).Value = "#string"
$ws.Cells.Item(2, 9).Value = "#float,  unit:µlormg"

# Row 3: new French descriptions row
$ws.Cells.Item(3, 1).Value = "#Manipulateur"
$ws.Cells.Item(3, 2).Value = "#Desc:IdentifiantEchantillon"
$ws.Cells.Item(3, 3).Value = "#Date"
$ws.Cells.Item(3, 4).Value = "#ModeOderatoireLaboratoire"
$ws.Cells.Item(3, 5).Value = "#AppareilLogicielCritique"
$ws.Cells.Item(3, 6).Value = "#ProduitCritique"
$ws.Cells.Item(3, 7).Value = "#LieuStockageDonneesBrutes"
$ws.Cells.Item(3, 8).Value = "#TypeExtraction"
$ws.Cells.Item(3, 9).Value = "#PriseEssai"
